$wb = $excel.ActiveWorkbook

$wsDual = $wb.Worksheets.Item(5)
$wsPlan = $wb.Worksheets.Item(6)

# ----------------------------------------------------------------------
# Sheet "XMOS Dualchip Planning" (6th sheet) -- add a "XSYS UART" row to
# the 1-bit port usage tables, pushing the totals rows down.
# ----------------------------------------------------------------------

# Insert a fresh row above the old "totals" row (row 13) to make room for
# the new "XSYS UART" usage entry; everything below shifts down by one.
$wsPlan.Rows.Item(13).Insert()

# New row 13: XSYS UART uses 2 one-bit ports (by block and by core).
$wsPlan.Range("A13").Value = "XSYS UART"
$wsPlan.Range("B13").Value = 2
$wsPlan.Range("G13").Value = "XSYS UART"
$wsPlan.Range("H13").Value = 2

# ----------------------------------------------------------------------
# Sheet "XMOS Dualchip" (5th sheet) -- JP1 connector now wired up like
# the XK-1A's UEXT/XSYS connector: add X0 Signal (column D) entries for
# P1E-P1J (rows 14-41 range of the "X0 Port" block).
# ----------------------------------------------------------------------

$wsDual.Range("D14").Value = "JP3_9"
$wsDual.Range("D15").Value = "JP3_7"
$wsDual.Range("D24").Value = "JP3_8"
$wsDual.Range("D25").Value = "XSYS_UART_TX"
$wsDual.Range("D26").Value = "JP3_10"
$wsDual.Range("D27").Value = "XSYS_UART_RX"
$wsDual.Range("D38").Value = "JP3_4"
$wsDual.Range("D39").Value = "JP3_3"
$wsDual.Range("D40").Value = "JP3_2"
$wsDual.Range("D41").Value = "JP3_1"

# Back to the planning sheet: row 14 now holds only the "by block" port
# total, extended to include the new row 13.
$wsPlan.Range("B14").Formula = "=SUM(B2:B13)"

# The "by core" block (previously sharing row 14 with the totals above)
# moves down to row 15, with its ranges extended to include row 13, and
# no longer uses a shared formula for I/J.
$wsPlan.Range("B15").ClearContents()
$wsPlan.Range("H15").Formula = "=SUM(H2:H13)"
$wsPlan.Range("I15").Formula = "=SUM(I2:I13)"
$wsPlan.Range("J15").Formula = "=SUM(J2:J13)"
$wsPlan.Range("K15").Formula = "=SUM(K2:K13)"
$wsPlan.Range("L15").Formula = "=SUM(L2:L13)"

# Row 16 ("1-bit ports remaining") already shifted down correctly with
# references pointing at row 15; nothing else to fix there.

# After the shift, row 17 is empty and row 18 holds both the
# "By Block ^" and "By Core^" labels (previously on row 17). Move "By
# Block ^" back up to the now-empty row 17, leaving "By Core^" alone on
# row 18.
$wsPlan.Range("A17").Value = "By Block ^"
$wsPlan.Range("A18").ClearContents()

# ----------------------------------------------------------------------
# Restore the view selections the way the author left them: last
# clicked around on the planning sheet, then returned focus to the
# dualchip connector sheet.
# ----------------------------------------------------------------------
$wsPlan.Activate()
$wsPlan.Range("L15").Select()

$wsDual.Activate()
$wsDual.Range("D28").Select()
